$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scale the "value" column (D) from units of 10k (wan) up by a factor of
# 10000, for every data row that actually has a numeric value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current * 10000
    }
}
